$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Phys_Stats_Table")

# --- Insert the new "Sheet1" worksheet, positioned after Phys_Stats_Table and before summary ---
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "Sheet1"

# Row map: destination row -> source row on Phys_Stats_Table
$rowMap = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 8
    6  = 9
    7  = 10
    8  = 11
    9  = 12
    10 = 13
    11 = 14
    12 = 15
    13 = 16
    14 = 17
    15 = 18
    16 = 19
    17 = 20
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    # columns B:C copy straight across
    $src.Range("B$srcRow`:C$srcRow").Copy($new.Range("B$destRow"))
    # three stat-group blocks are reordered: source X:AA, P:S, T:W -> dest D:G, H:K, L:O
    $src.Range("X$srcRow`:AA$srcRow").Copy($new.Range("D$destRow"))
    $src.Range("P$srcRow`:S$srcRow").Copy($new.Range("H$destRow"))
    $src.Range("T$srcRow`:W$srcRow").Copy($new.Range("L$destRow"))
}

# Footnote row
$src.Range("B72:C72").Copy($new.Range("B18"))

# "Two-way ANOVA" header, re-typed using the style already used for the matching
# "Three-way ANOVA" header cell (P12 / dest H9)
$new.Range("H4").Value = $src.Range("O7").Value
$src.Range("P12").Copy($new.Range("H4"))
$new.Range("H4").Value = $src.Range("O7").Value

Write-Output "done"
